$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 24 (2025-11 stats) with new values
$ws.Range("B24").Value = 6354
$ws.Range("C24").Value = 998
$ws.Range("D24").Value = 5955260
$ws.Range("E24").Value = 937.2458293988039
$ws.Range("F24").Value = 8.319127173542441
$ws.Range("G24").Value = 3.419689119170988
$ws.Range("H24").Value = 26.15547437895238
